# "add extended test cases" — append new meeting-notes entries and
# re-prioritise a handful of existing TODO rows on the "Meeting Notes" sheet;
# also nudge the view/selection on "Meeting Notes" and "Tests".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Meeting Notes")

# --- Re-order the existing open-questions list (rows 41-45, column C) ---
$ws.Range("C41").Value = "unnesting 100% (multiple correlations,all subquery operators, ...)"
$ws.Range("C42").Value = "add schema"
$ws.Range("C43").Value = "clean up code (c-style pointers, error handling)"
$ws.Range("C44").Value = "build UI"
$ws.Range("C45").Value = "visualize RA tree"

# --- Append the new meeting entry (rows 48-53) ---
# Cells are written in the same order the shared-string table records them
# in, so new unique strings land at the expected indices.
$ws.Range("A48").Value = "20.12.2022"
$ws.Range("B53").Value = "MA"
$ws.Range("B49").Value = "Tests to 100% algorithm"
$ws.Range("C49").Value = "outer/semi joins"
$ws.Range("C50").Value = "multiple correlations in query"
$ws.Range("B48").Value = "Q1/Q2 & TPCH decorrelated (special handling for exists/in)"
$ws.Range("C51").Value = "set operations (add to sql parser)"
$ws.Range("B52").Value = "Q1 läuft in Umbra nicht (siehe screenshots)"
$ws.Range("C52").Value = "mit CTEs geht's schon"

# --- View state: scroll + selection on "Meeting Notes" ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 33
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C47").Select()

# --- View state: selection on "Tests" ---
$ws2 = $wb.Worksheets.Item("Tests")
$ws2.Activate()
$ws2.Range("D10").Select()

# Leave "Meeting Notes" as the active/visible sheet, matching tabSelected.
$ws.Activate()
